$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 12499
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 12499
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 37497
$ws.Range("N69").Value = -39245

$ws.Range("H70").Value = 11831.8
$ws.Range("I70").Value = 2310
$ws.Range("J70").Value = 26114.5
$ws.Range("K70").Value = 6930
$ws.Range("L70").Value = 78343.5
$ws.Range("M70").Value = -6660
$ws.Range("N70").Value = -78883.5

$ws.Range("H72").Value = 12499
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 12499
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 112491
$ws.Range("N72").Value = -121227

$ws.Range("H73").Value = 11831.8
$ws.Range("I73").Value = 2310
$ws.Range("J73").Value = 26114.5
$ws.Range("K73").Value = 6930
$ws.Range("L73").Value = 78343.5
$ws.Range("M73").Value = -5994
$ws.Range("N73").Value = -80215.5

$ws.Range("H76").Value = 9245.454
$ws.Range("I76").Value = 9111.5
$ws.Range("J76").Value = 9602.666999999999
$ws.Range("K76").Value = 9111.5
$ws.Range("L76").Value = 9602.666999999999
$ws.Range("M76").Value = -8796.5
$ws.Range("N76").Value = -10232.667

$ws.Range("H79").Value = 9245.454
$ws.Range("I79").Value = 9111.5
$ws.Range("J79").Value = 9602.666999999999
$ws.Range("K79").Value = 9111.5
$ws.Range("L79").Value = 9602.666999999999
$ws.Range("M79").Value = -8019.5
$ws.Range("N79").Value = -11786.667

$ws.Range("H88").Value = 1973.7646
$ws.Range("I88").Value = 2781.2856
$ws.Range("J88").Value = 1408.5
$ws.Range("K88").Value = 2781.2856
$ws.Range("L88").Value = 1408.5
$ws.Range("M88").Value = -2375.2856
$ws.Range("N88").Value = -2220.5

$ws.Range("H91").Value = 1973.7646
$ws.Range("I91").Value = 2781.2856
$ws.Range("J91").Value = 1408.5
$ws.Range("K91").Value = 2781.2856
$ws.Range("L91").Value = 1408.5
$ws.Range("M91").Value = -1377.2856
$ws.Range("N91").Value = -4216.5

$ws.Range("H101").Value = 2488.4
$ws.Range("I101").Value = 2839.25
$ws.Range("J101").Value = 1085
$ws.Range("K101").Value = 8517.75
$ws.Range("L101").Value = 3255
$ws.Range("M101").Value = -6895.75
$ws.Range("N101").Value = -6499

$ws.Range("H106").Value = 29480.428
$ws.Range("I106").Value = 29480.428
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 29480.428
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -28849.428

$ws.Range("H111").Value = 1696.1428
$ws.Range("I111").Value = 899
$ws.Range("J111").Value = 1829
$ws.Range("K111").Value = 2697
$ws.Range("L111").Value = 5487
$ws.Range("M111").Value = 370
$ws.Range("N111").Value = -11621

$ws.Range("H112").Value = 5247.2266
$ws.Range("I112").Value = 14977.8
$ws.Range("J112").Value = 4233.625
$ws.Range("K112").Value = 44933.39999999999
$ws.Range("L112").Value = 12700.875
$ws.Range("M112").Value = -43825.39999999999
$ws.Range("N112").Value = -14916.875

$ws.Range("H138").Value = 11357.315
$ws.Range("I138").Value = 11499.667
$ws.Range("J138").Value = 11291.615
$ws.Range("K138").Value = 34499.001
$ws.Range("L138").Value = 33874.845
$ws.Range("M138").Value = -29359.001
$ws.Range("N138").Value = -44154.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2158840.5
$ws.Range("I32").Value = 3237758.5
$ws.Range("J32").Value = 68436.625
$ws.Range("K32").Value = 3237758.5
$ws.Range("L32").Value = 68436.625
$ws.Range("M32").Value = -3237471.5
$ws.Range("N32").Value = -69010.625

$ws.Range("H61").Value = 15450.447
$ws.Range("I61").Value = 4176.8887
$ws.Range("J61").Value = 43121.91
$ws.Range("K61").Value = 4176.8887
$ws.Range("L61").Value = 43121.91
$ws.Range("M61").Value = -3964.8887
$ws.Range("N61").Value = -43545.91

$ws.Range("H74").Value = 19145.72
$ws.Range("I74").Value = 2677.4546
$ws.Range("J74").Value = 32085.072
$ws.Range("K74").Value = 2677.4546
$ws.Range("L74").Value = 32085.072
$ws.Range("M74").Value = -1803.4546
$ws.Range("N74").Value = -33833.072

$ws.Range("H77").Value = 19145.72
$ws.Range("I77").Value = 2677.4546
$ws.Range("J77").Value = 32085.072
$ws.Range("K77").Value = 13387.273
$ws.Range("L77").Value = 160425.36
$ws.Range("M77").Value = -9019.273000000001
$ws.Range("N77").Value = -169161.36

$ws.Range("H110").Value = 3515.111
$ws.Range("I110").Value = 1771.25
$ws.Range("J110").Value = 17466
$ws.Range("K110").Value = 1771.25
$ws.Range("L110").Value = 17466
$ws.Range("M110").Value = 273.75
$ws.Range("N110").Value = -21556

$ws.Range("H136").Value = 15450.447
$ws.Range("I136").Value = 4176.8887
$ws.Range("J136").Value = 43121.91
$ws.Range("K136").Value = 12530.6661
$ws.Range("L136").Value = 129365.73
$ws.Range("M136").Value = -9980.666100000002
$ws.Range("N136").Value = -134465.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 375
$ws.Range("I18").Value = 350
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 350
$ws.Range("L18").Value = 400
$ws.Range("M18").Value = 179
$ws.Range("N18").Value = -1458

$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 1250
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1250
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1077

$ws.Range("H26").Value = 26788.2
$ws.Range("I26").Value = 26788.2
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 26788.2
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -26496.2

$ws.Range("H96").Value = 32206
$ws.Range("I96").Value = 32206
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 32206
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -29460

$ws.Range("H115").Value = 160000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 160000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 160000
$ws.Range("N115").Value = -163134

$ws.Range("H134").Value = 10747.021
$ws.Range("I134").Value = 6653.5557
$ws.Range("J134").Value = 24143.818
$ws.Range("K134").Value = 19960.6671
$ws.Range("L134").Value = 72431.454
$ws.Range("M134").Value = -17425.6671
$ws.Range("N134").Value = -77501.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 442.65
$ws.Range("I7").Value = 103.4
$ws.Range("J7").Value = 1460.4
$ws.Range("K7").Value = 103.4
$ws.Range("L7").Value = 1460.4
$ws.Range("M7").Value = 9.599999999999994
$ws.Range("N7").Value = -1686.4

$ws.Range("H31").Value = 25521.908
$ws.Range("I31").Value = 16824
$ws.Range("J31").Value = 35048.19
$ws.Range("K31").Value = 16824
$ws.Range("L31").Value = 35048.19
$ws.Range("M31").Value = -16529
$ws.Range("N31").Value = -35638.19

$ws.Range("H34").Value = 25521.908
$ws.Range("I34").Value = 16824
$ws.Range("J34").Value = 35048.19
$ws.Range("K34").Value = 16824
$ws.Range("L34").Value = 35048.19
$ws.Range("M34").Value = -16622
$ws.Range("N34").Value = -35452.19

$ws.Range("H58").Value = 10012.745
$ws.Range("I58").Value = 5161.96
$ws.Range("J58").Value = 15525
$ws.Range("K58").Value = 5161.96
$ws.Range("L58").Value = 15525
$ws.Range("M58").Value = -4958.96
$ws.Range("N58").Value = -15931

$ws.Range("H105").Value = 14542.917
$ws.Range("I105").Value = 16065.143
$ws.Range("J105").Value = 12411.8
$ws.Range("K105").Value = 16065.143
$ws.Range("L105").Value = 12411.8
$ws.Range("M105").Value = -14318.143
$ws.Range("N105").Value = -15905.8

$ws.Range("H132").Value = 6920.1387
$ws.Range("I132").Value = 1666.8
$ws.Range("J132").Value = 13486.8125
$ws.Range("K132").Value = 5000.4
$ws.Range("L132").Value = 40460.4375
$ws.Range("M132").Value = -2470.4
$ws.Range("N132").Value = -45520.4375

$ws.Range("H136").Value = 10012.745
$ws.Range("I136").Value = 5161.96
$ws.Range("J136").Value = 15525
$ws.Range("K136").Value = 15485.88
$ws.Range("L136").Value = 46575
$ws.Range("M136").Value = -12935.88
$ws.Range("N136").Value = -51675

$ws.Range("H141").Value = 111317.43
$ws.Range("I141").Value = 99890.5
$ws.Range("J141").Value = 115888.2
$ws.Range("K141").Value = 99890.5
$ws.Range("L141").Value = 115888.2
$ws.Range("M141").Value = -94710.5
$ws.Range("N141").Value = -126248.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1995
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5985
$ws.Range("N113").Value = -10325

$ws.Range("H122").Value = 23922612
$ws.Range("I122").Value = 46717670
$ws.Range("J122").Value = 5686564
$ws.Range("K122").Value = 420459030
$ws.Range("L122").Value = 51179076
$ws.Range("M122").Value = -420456580
$ws.Range("N122").Value = -51183976

$ws.Range("H129").Value = 9094423
$ws.Range("I129").Value = 33334250
$ws.Range("J129").Value = 4487.875
$ws.Range("K129").Value = 100002750
$ws.Range("L129").Value = 13463.625
$ws.Range("M129").Value = -99997750
$ws.Range("N129").Value = -23463.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13054.857
$ws.Range("I126").Value = 3698.5
$ws.Range("J126").Value = 16797.4
$ws.Range("K126").Value = 11095.5
$ws.Range("L126").Value = 50392.2
$ws.Range("M126").Value = -8625.5
$ws.Range("N126").Value = -55332.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9169.666999999999
$ws.Range("I22").Value = 5393
$ws.Range("J22").Value = 20499.666
$ws.Range("K22").Value = 5393
$ws.Range("L22").Value = 20499.666
$ws.Range("M22").Value = -5098
$ws.Range("N22").Value = -21089.666

$ws.Range("H27").Value = 9169.666999999999
$ws.Range("I27").Value = 5393
$ws.Range("J27").Value = 20499.666
$ws.Range("K27").Value = 5393
$ws.Range("L27").Value = 20499.666
$ws.Range("M27").Value = -5286
$ws.Range("N27").Value = -20713.666

$ws.Range("H100").Value = 16333.167
$ws.Range("I100").Value = 30000
$ws.Range("J100").Value = 13599.8
$ws.Range("K100").Value = 30000
$ws.Range("L100").Value = 13599.8
$ws.Range("M100").Value = -29459
$ws.Range("N100").Value = -14681.8

$ws.Range("H120").Value = 139666.33
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 139666.33
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 139666.33
$ws.Range("N120").Value = -149342.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H52").Value = 19598.4
$ws.Range("I52").Value = 10997.333
$ws.Range("J52").Value = 32500
$ws.Range("K52").Value = 10997.333
$ws.Range("L52").Value = 32500
$ws.Range("M52").Value = -10771.333
$ws.Range("N52").Value = -32952

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 117242.375
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 117242.375
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 117242.375
$ws.Range("N64").Value = -117738.375

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 117242.375
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 117242.375
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 117242.375
$ws.Range("N67").Value = -118958.375

$ws.Range("H132").Value = 6150.636
$ws.Range("I132").Value = 2688.375
$ws.Range("J132").Value = 15383.333
$ws.Range("K132").Value = 8065.125
$ws.Range("L132").Value = 46149.999
$ws.Range("M132").Value = -5535.125
$ws.Range("N132").Value = -51209.999
